{"js": "// Office.js (Word JavaScript API) script.\n// Adds the \"Plans and Progress\" planning/progress paragraphs that were\n// inserted just before the existing \"Frecipe will require several key\n// components;\" paragraph, and relocates the Word-managed \"_GoBack\"\n// bookmark to sit immediately before the \"Project management\" bullet\n// (its new position after the edit), matching the author's edit.\n\n// 1) Locate the existing paragraph that starts the \"key components\" list\n//    intro sentence - the anchor for our insertion.\nconst anchorResults = context.document.body.search(\n  \"Frecipe will require several key components;\",\n  { matchCase: true }\n);\nanchorResults.load(\"text\");\nawait context.sync();\n\nif (anchorResults.items.length === 0) {\n  throw new Error('Could not find anchor paragraph \"Frecipe will require several key components;\"');\n}\n\nconst anchorParagraph = anchorResults.items[0].paragraphs.getFirst();\n\n// 2) Insert the three new, fully self-contained paragraphs directly\n//    before the anchor paragraph (in document order).\nconst paraIntro = anchorParagraph.insertParagraph(\n  \"Frecipe began as an aid to deciding what meal to prepare based on what food was in your fridge\",\n  Word.InsertLocation.before\n);\nawait context.sync();\n\n// Second sentence of the same paragraph (kept as its own insertText call\n// so it is appended after the first sentence rather than merged by hand).\nparaIntro.insertText(\n  \". Along the planning process we discovered the idea has a lot more potential than just helping people decide what meal to prepare. \",\n  Word.InsertLocation.end\n);\nawait context.sync();\n\nanchorParagraph.insertParagraph(\n  \"We discovered during our planning stage the idea has potential to expand outside of the fridge and incorporate the pantry into the picture. \",\n  Word.InsertLocation.before\n);\nawait context.sync();\n\nanchorParagraph.insertParagraph(\n  \"Most importantly, we established the app has the potential to make a big impact in reducing food waste as we forecast if people use the app they are more likely to use items in their fridge and pantry rather than eat out all the time which also in turn saves  them money. \",\n  Word.InsertLocation.before\n);\nawait context.sync();\n\n// 3) Prepend \"During planning we established \" to the existing anchor\n//    paragraph text, so it reads \"During planning we established Frecipe\n//    will require several key components;\" - same paragraph as before.\n//    NOTE: Paragraph.insertText only accepts Start/End/Replace locations;\n//    Start behaves as an in-place prepend (unlike Range.insertText, which\n//    also supports Before/After).\nanchorParagraph.insertText(\"During planning we established \", Word.InsertLocation.start);\nawait context.sync();\n\n// 4) Move the \"_GoBack\" bookmark (Word's \"last edit location\" marker) to\n//    sit right before the \"Project management\" bullet item, which is\n//    where it ends up after this edit.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst pmResults = context.document.body.search(\"Project management\", { matchCase: true });\npmResults.load(\"text\");\nawait context.sync();\n\nif (pmResults.items.length > 0) {\n  const pmParagraph = pmResults.items[0].paragraphs.getFirst();\n  const pmStart = pmParagraph.getRange(Word.RangeLocation.start);\n  pmStart.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Adds the \"Plans and Progress\" planning/progress paragraphs that were\n# inserted just before the existing \"Frecipe will require several key\n# components;\" paragraph, and relocates the Word-managed \"_GoBack\"\n# bookmark to sit immediately before the \"Project management\" bullet\n# (its new position after the edit), matching the author's edit.\n\n$d = $word.ActiveDocument\n$cr = [char]13\n\n$anchorText = \"Frecipe will require several key components;\"\n\n# 1) Locate the anchor paragraph/sentence.\n$rng = $d.Content\n$rng.Find.Execute($anchorText) | Out-Null\n\n# 2) Insert the three new, self-contained paragraphs directly before the\n#    anchor paragraph (in document order). Re-finding the anchor text\n#    after each mutation keeps the range accurate.\n$para1 = \"Frecipe began as an aid to deciding what meal to prepare based on what food was in your fridge. Along the planning process we discovered the idea has a lot more potential than just helping people decide what meal to prepare. \"\n$rng.InsertBefore($para1 + $cr)\n\n$rng = $d.Content\n$rng.Find.Execute($anchorText) | Out-Null\n$para2 = \"We discovered during our planning stage the idea has potential to expand outside of the fridge and incorporate the pantry into the picture. \"\n$rng.InsertBefore($para2 + $cr)\n\n$rng = $d.Content\n$rng.Find.Execute($anchorText) | Out-Null\n$para3 = \"Most importantly, we established the app has the potential to make a big impact in reducing food waste as we forecast if people use the app they are more likely to use items in their fridge and pantry rather than eat out all the time which also in turn saves  them money. \"\n$rng.InsertBefore($para3 + $cr)\n\n# 3) Prepend \"During planning we established \" to the existing anchor\n#    paragraph text (no new paragraph mark), so it reads \"During planning\n#    we established Frecipe will require several key components;\".\n$rng = $d.Content\n$rng.Find.Execute($anchorText) | Out-Null\n$rng.InsertBefore(\"During planning we established \")\n\n# 4) Move the \"_GoBack\" bookmark (Word's \"last edit location\" marker) to\n#    sit right before the \"Project management\" bullet item, which is\n#    where it ends up after this edit.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$pmRng = $d.Content\n$pmRng.Find.Execute(\"Project management\") | Out-Null\n$pmTarget = $d.Range($pmRng.Start, $pmRng.Start)\n$d.Bookmarks.Add(\"_GoBack\", $pmTarget)\n"}
